# Apply the text replacements described by the diff.
# Word's Find/Replace (wdReplaceOne) replaces only the first match, which is
# safe here because every source string is unique at the time it is searched
# for (the one potential collision -- "55÷7=" becoming "14÷8=" while an
# original "14÷8=" cell also needs to change -- is avoided by replacing the
# original "14÷8=" cell *before* the "55÷7=" cell is turned into "14÷8=").

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2023-11-01 Wednesday" "2023-11-02 Thursday"

# Row 1
Replace-Text "66÷4=" "32÷8="
Replace-Text "88÷2=" "55÷4="
Replace-Text "14÷7=" "16÷7="
Replace-Text "49÷8=" "17÷5="
Replace-Text "59÷8=" "44÷9="

# Row 2
Replace-Text "25÷9=" "73÷4="
Replace-Text "86÷4=" "40÷3="
Replace-Text "65÷7=" "37÷2="
Replace-Text "50÷6=" "44÷3="
Replace-Text "98÷2=" "22÷4="

# Row 3
Replace-Text "23÷8=" "56÷4="
Replace-Text "35÷2=" "77÷5="
Replace-Text "27÷8=" "57÷6="
Replace-Text "99÷2=" "97÷9="
Replace-Text "30÷9=" "59÷4="

# Row 4 (note: "14÷8=" must be replaced before "55÷7=" below produces a new
# "14÷8=" text, otherwise the later search would erroneously match it)
Replace-Text "20÷3=" "90÷5="
Replace-Text "49÷9=" "61÷2="
Replace-Text "14÷8=" "10÷8="
Replace-Text "99÷3=" "74÷7="
Replace-Text "31÷2=" "38÷4="

# Row 5
Replace-Text "43÷3=" "33÷7="
Replace-Text "55÷7=" "14÷8="
Replace-Text "12÷8=" "96÷3="
Replace-Text "52÷2=" "73÷8="
Replace-Text "21÷7=" "16÷3="
